$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the regression coefficient values (row 2 = A Lag, row 3 = C Lag,
# row 4 = FFR Lag, row 5 = LF Lag; columns B..E = A, C, FFR, LF).
# Some replacement values look like plain numbers (e.g. "1.65"); prefix
# with an apostrophe so Excel keeps them as text (matching the original
# shared-string / text cell type), then reset the style so no visible
# "number stored as text" formatting marker is left behind.
$cells = [ordered]@{
    "B2" = "-0.29***"
    "B3" = "-1.02*"
    "B4" = "0.06***"
    "B5" = "0.01***"
    "C2" = "-0.02***"
    "C3" = "-0.42***"
    "C4" = "-0.0*"
    "C5" = "0.0***"
    "D2" = "0.26*"
    "D3" = "1.65"
    "D4" = "0.32***"
    "D5" = "-0.01*"
    "E2" = "-3.7*"
    "E3" = "-1.49"
    "E4" = "0.53"
    "E5" = "0.1"
}

foreach ($addr in $cells.Keys) {
    $ws.Range($addr).Formula = "'" + $cells[$addr]
    $ws.Range($addr).Style = "Normal"
}
